{"js": "// 1) Insert a new centered \"\u0412\u0430\u0440\u0456\u0430\u043d\u0442 \u211619\" paragraph right after the \"\u0422\u0435\u043c\u0430 ...\"\n//    heading paragraph (i.e. immediately before the blank centered paragraph\n//    that currently follows it).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet temaParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"\u0421\u0442\u0430\u0442\u0438\u0441\u0442\u0438\u0447\u043d\u0438\u0439 \u0430\u043d\u0430\u043b\u0456\u0437 \u0456 \u043f\u0435\u0440\u0432\u0438\u043d\u043d\u0430 \u043e\u0431\u0440\u043e\u0431\u043a\u0430 \u0434\u0430\u043d\u0438\u0445\") !== -1) {\n    temaParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!temaParagraph) {\n  throw new Error(\"Could not locate the '\u0422\u0435\u043c\u0430 ...' heading paragraph\");\n}\n\nconst variantParagraph = temaParagraph.insertParagraph(\"\u0412\u0430\u0440\u0456\u0430\u043d\u0442 \u211619\", Word.InsertLocation.after);\nvariantParagraph.alignment = Word.Alignment.centered;\nvariantParagraph.firstLineIndent = 0;\nawait context.sync();\n\n// 2) Remove the empty \"P\"-styled paragraph that directly follows the \"2024\"\n//    paragraph (collapsing it into the surrounding content, matching the\n//    removed <w:p> block in the target diff).\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text,items/style\");\nawait context.sync();\n\nlet emptyPParagraph = null;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].style === \"P\" && paragraphs2.items[i].text === \"\") {\n    emptyPParagraph = paragraphs2.items[i];\n    break;\n  }\n}\nif (!emptyPParagraph) {\n  throw new Error(\"Could not locate the empty 'P'-styled paragraph to remove\");\n}\n\nemptyPParagraph.delete();\nawait context.sync();\n", "ps1": "# 1) Insert a new centered \"\u0412\u0430\u0440\u0456\u0430\u043d\u0442 \u211619\" paragraph right after the \"\u0422\u0435\u043c\u0430 ...\"\n#    heading paragraph (i.e. immediately before the blank centered paragraph\n#    that currently follows it).\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$temaIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*\u0421\u0442\u0430\u0442\u0438\u0441\u0442\u0438\u0447\u043d\u0438\u0439 \u0430\u043d\u0430\u043b\u0456\u0437 \u0456 \u043f\u0435\u0440\u0432\u0438\u043d\u043d\u0430 \u043e\u0431\u0440\u043e\u0431\u043a\u0430 \u0434\u0430\u043d\u0438\u0445*\") {\n        $temaIndex = $i\n        break\n    }\n}\nif ($temaIndex -eq -1) {\n    throw \"Could not locate the '\u0422\u0435\u043c\u0430 ...' heading paragraph\"\n}\n\n$temaPara = $d.Paragraphs.Item($temaIndex)\n$temaPara.Range.InsertParagraphAfter()\n\n$variantPara = $d.Paragraphs.Item($temaIndex + 1)\n$variantPara.Range.Text = \"\u0412\u0430\u0440\u0456\u0430\u043d\u0442 \u211619\"\n$variantPara.Alignment = 1\n\n# 2) Remove the empty \"P\"-styled paragraph that directly follows the \"2024\"\n#    paragraph (collapsing it into the surrounding content, matching the\n#    removed <w:p> block in the target diff).\n$count2 = $d.Paragraphs.Count\n$yearIndex = -1\nfor ($i = 1; $i -le $count2; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"2024*\") {\n        $yearIndex = $i\n        break\n    }\n}\nif ($yearIndex -eq -1) {\n    throw \"Could not locate the '2024' paragraph\"\n}\n\n$nextPara = $d.Paragraphs.Item($yearIndex + 1)\nif ($nextPara.Style.NameLocal -eq \"P\") {\n    $nextPara.Range.Delete()\n} else {\n    throw \"Paragraph following '2024' was not the expected empty 'P'-styled paragraph\"\n}\n"}
